$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SignUpPage")

# Add new row 8, mirroring the pattern of row 7 but with a new unique label in column A.
$ws.Range("A8").Value = "MeasurabilityFlow3"
$ws.Range("A8").NumberFormat = "@"

$ws.Range("B8").Style = "Normal"
$ws.Range("B8").Value = "embibe.auto2@mailinator.com"

$ws.Range("C8").Value = "embibe123"
$ws.Range("C8").NumberFormat = "@"

$ws.Range("D8").Value = "Engineering"
$ws.Range("D8").NumberFormat = "@"

# Move the active selection to B5, matching the saved workbook state.
$ws.Range("B5").Select()
